$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.718.24"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "2.301.27"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'300.96"
$ws.Range("E5").Value = "  -1.54%  "

$ws.Range("D6").Value = "'96.05"
$ws.Range("E6").Value = "  -1.51%  "

$ws.Range("D7").Value = "'0.504"
$ws.Range("E7").Value = "  -1.46%  "

$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").Value = "'0.496"
$ws.Range("E9").Value = "  -1.87%  "

$ws.Range("D10").Value = "'34.64"
$ws.Range("E10").Value = "  -3.18%  "

$ws.Range("D11").Value = "'19.22"
$ws.Range("E11").Value = "  +5.68%  "

$ws.Range("D12").Value = "'0.0788"
$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").Value = "'6.81"
$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "2.668.02"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").Value = "2.312.19"
$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").Value = "'0.782"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "42.749.48"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").Value = "'12.34"
$ws.Range("E19").Value = "  -2.95%  "

$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").Value = "  -1.87%  "

$ws.Range("D21").Value = "'6.01"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").Value = "'67.46"
$ws.Range("E22").Value = "  -0.48%  "

$ws.Range("D23").Value = "'235.13"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  +2.94%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  -2.06%  "

$ws.Range("D27").Value = "'24.59"
$ws.Range("E27").Value = "  -3.47%  "

$ws.Range("D28").Value = "'2.20"
$ws.Range("E28").Value = "  +6.44%  "

$ws.Range("D29").Value = "'163.88"
$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("D30").Value = "'9.03"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").Value = "'32.34"
$ws.Range("E31").Value = "  -2.33%  "

$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").Value = "'4.95"
$ws.Range("E33").Value = "  -1.33%  "

$ws.Range("D34").Value = "'17.60"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").Value = "'4.47"
$ws.Range("E35").Value = "  -7.45%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.34"
$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.0695"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("E38").Value = "  -1.17%  "

$ws.Range("D39").Value = "'1.75"
$ws.Range("E39").Value = "  -0.41%  "

$ws.Range("D40").Value = "'2.74"
$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("D41").Value = "'0.108"
$ws.Range("E41").Value = "  -1.31%  "

$ws.Range("D42").Value = "1.971.29"
$ws.Range("E42").Value = "  -1.83%  "

$ws.Range("D43").Value = "'10.53"
$ws.Range("E43").Value = "  +5.08%  "

$ws.Range("D44").Value = "'18.62"
$ws.Range("E44").Value = "  +3.77%  "

$ws.Range("D45").Value = "'0.0278"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("E46").Value = "  -2.62%  "

$ws.Range("D47").Value = "'2.75"
$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("D48").Value = "'2.91"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").Value = "2.535.66"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").Value = "'53.02"
$ws.Range("E50").Value = "  -2.18%  "

$ws.Range("D51").Value = "'72.11"
$ws.Range("E51").Value = "  +0.12%  "
